$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.087.20'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +4.78%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.240.96'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +4.88%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '250.99'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +7.21%  '
$ws.Range('E6').Value = '  +2.73%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '75.26'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +9.78%  '
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('E9').Value = '  +6.50%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.17'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +8.25%  '
$ws.Range('E11').Value = '  +4.78%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.91'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +5.56%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.101'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.17%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.576.18'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.80%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.65'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.01%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.239.01'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +4.56%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.792'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.71%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '42.970.91'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +4.83%  '
$ws.Range('E19').Value = '  +6.71%  '
$ws.Range('E20').Value = '  +3.21%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.99'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +5.44%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '229.58'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.89%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.19'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +17.27%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.65'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.51%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.78'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.91%  '
$ws.Range('E27').Value = '  +2.65%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '39.42'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +28.52%  '
$ws.Range('E29').Value = '  +5.94%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.17'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.26%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '171.42'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '20.22'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +4.11%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0800'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +7.41%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.26'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.07%  '
$ws.Range('E35').Value = '  +2.53%  '
$ws.Range('E36').Value = '  +11.06%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.46'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +10.22%  '
$ws.Range('E38').Value = '  +18.41%  '
$ws.Range('E39').Value = '  +13.02%  '
$ws.Range('E40').Value = '  +4.72%  '
$ws.Range('E41').Value = '  +11.35%  '
$ws.Range('E42').Value = '  +4.28%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '59.48'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +5.06%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.66'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +6.56%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '103.81'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +8.06%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.480'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +30.81%  '
$ws.Range('E47').Value = '  +4.59%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.43'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +14.63%  '
$ws.Range('E49').Value = '  +4.42%  '
$ws.Range('E50').Value = '  +5.41%  '
$ws.Range('E51').Value = '  +3.30%  '
